$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty "Survey 3" data row (row 4)
$ws.Range("B4").Value = 13
$ws.Range("C4").Value = 40
$ws.Range("D4").Value = 22
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0

# Update the active selection to match the saved workbook state
$ws.Range("D8").Select() | Out-Null
